# Update of the Slovakia Covid daily stats sheet ("pi 06. 05. 2022")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing AgTests (F) / AgPosit (G) figures for prior days ---
$updates = @(
    @{ Row = 699; F = 43675 },
    @{ Row = 714; F = 32760 },
    @{ Row = 715; F = 32073 },
    @{ Row = 716; F = 30007 },
    @{ Row = 719; F = 45312 },
    @{ Row = 758; F = 11225; G = 931 },
    @{ Row = 759; F = 3875 },
    @{ Row = 760; F = 5162; G = 565 },
    @{ Row = 761; F = 16986; G = 1286 },
    @{ Row = 762; F = 11316; G = 799 },
    @{ Row = 763; F = 10366; G = 735 },
    @{ Row = 764; F = 11193; G = 695 },
    @{ Row = 765; F = 9222; G = 610 },
    @{ Row = 767; F = 4145 },
    @{ Row = 768; F = 15072; G = 778 },
    @{ Row = 769; F = 10106; G = 515 },
    @{ Row = 770; F = 9153; G = 429 },
    @{ Row = 771; F = 9240; G = 410 },
    @{ Row = 776; F = 14887; G = 674 },
    @{ Row = 777; F = 10441; G = 453 },
    @{ Row = 778; F = 9050; G = 360 },
    @{ Row = 779; F = 7376; G = 305 },
    @{ Row = 780; F = 2694 },
    @{ Row = 781; F = 2750; G = 152 },
    @{ Row = 782; F = 10715; G = 427 },
    @{ Row = 783; F = 7810; G = 249 },
    @{ Row = 784; F = 7724; G = 235 },
    @{ Row = 785; F = 7141; G = 236 },
    @{ Row = 786; F = 6284; G = 192 },
    @{ Row = 787; F = 2251; G = 94 },
    @{ Row = 788; F = 1748; G = 84 },
    @{ Row = 789; F = 7671; G = 325 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    if ($u.ContainsKey("G")) {
        $ws.Cells.Item($u.Row, 7).Value = $u.G
    }
}

# --- Fill in the previously-missing AgTests/AgPosit values for row 790 ---
$ws.Cells.Item(790, 6).Value = 4484
$ws.Cells.Item(790, 7).Value = 164

# --- Append two new daily rows (06.05.2022 commit adds 04.05 and 05.05.2022 data) ---
$newRows = @(
    @{ Row = 791; A = 44685; B = 1783906; C = 3051; D = 553;  E = 19952; F = 2795; G = 148 },
    @{ Row = 792; A = 44686; B = 1784410; C = 2436; D = 504;  E = 19961; F = 48;   G = 0 }
)

foreach ($nr in $newRows) {
    $ws.Cells.Item($nr.Row, 1).Value = $nr.A
    $ws.Cells.Item($nr.Row, 1).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($nr.Row, 2).Value = $nr.B
    $ws.Cells.Item($nr.Row, 3).Value = $nr.C
    $ws.Cells.Item($nr.Row, 4).Value = $nr.D
    $ws.Cells.Item($nr.Row, 5).Value = $nr.E
    $ws.Cells.Item($nr.Row, 6).Value = $nr.F
    $ws.Cells.Item($nr.Row, 7).Value = $nr.G
}
